$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Yes" is the first brand-new shared string the author introduced (it ends
# up earliest in the shared-string table), so stamp the Bonus column first.
$ws.Range("E15").Value = "Yes"
$ws.Range("E16").Value = "Yes"

# --- Row 15: "Added Wilsons algorithm" entry ---
$ws.Range("A15").Value = "Added Wilsons algorithm"
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = 45439
$ws.Range("D15").Value = "Wanted to really challenge myself so I looked up what  the hardest `nalgorithm to implement was according to others online. They all seemed`nto agree that Wilsons was the hardest to get right so I figured I'll add that`none. They were not lying as it's quite a tricky one but I managed to `nimplement it nonetheless."
$ws.Range("D15").WrapText = $true

# --- Row 16: "Added skybox and gradient" entry ---
$ws.Range("A16").Value = "Added skybox and gradient"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 45439
$ws.Range("D16").Value = "Added a skybox shader I tend to re-use a lot for small projects and used`nDoTween to attach a animated gradient to all the wall materials."
$ws.Range("D16").WrapText = $true

# Match the autofit row heights Excel would have computed for the newly
# wrapped, multi-line content.
$ws.Rows(15).RowHeight = 90.75
$ws.Rows(16).RowHeight = 27

# --- View state: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H15").Select()
